# Corrected typo in ACS slide:
#   - Slide 12: lower-case the word "Ubuntu" -> "ubuntu" inside the
#     'Pull "Ubuntu" image from Docker Hub or local registry' caption.
#   - Slide 16: tidy up the ssh command line so the port-forwarding
#     argument " -p 2200 -L 22375:127.0.0.1:2375" is one contiguous run
#     again instead of being split across two runs.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 - "Running a Container" - TextBox 14 (shape 13)
# ---------------------------------------------------------------------
$s12  = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(13)
$tr12 = $sh12.TextFrame.TextRange

# Original text: Pull "Ubuntu" image from Docker Hub or local registry
# Replace the capitalised "Ubuntu" (characters 7-12) with lower-case "ubuntu"
$tr12.Characters(7, 6).Text = "ubuntu"

# Re-assign the opening quote (character 6) on its own so that it becomes
# its own run, separate from the leading "Pull " text.
$tr12.Characters(6, 1).Text = [string][char]34

# ---------------------------------------------------------------------
# Slide 16 - "Connecting to Docker Swarm in ACS" - TextBox 3 (shape 3)
# ---------------------------------------------------------------------
$s16  = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(3)
$tr16 = $sh16.TextFrame.TextRange

# This shape auto-fits its height to the single line of text it contains.
# Remember the current height so we can put it back exactly afterwards -
# merging the runs below shouldn't visually change the line, so the
# shape's size must not change either.
$origHeight16 = $sh16.Height

# Suspend "resize shape to fit text" while we touch the runs so the
# in-progress edit can't leave the box in a half-updated state.
$sh16.TextFrame.AutoSize = 0

$full16 = $tr16.Text
$idx16  = $full16.IndexOf(" -p 2200 -L ")
# Re-write " -p 2200 -L " + "22375:127.0.0.1:2375" as a single run since
# both pieces already share identical formatting.
$tr16.Characters($idx16 + 1, 32).Text = " -p 2200 -L 22375:127.0.0.1:2375"

# Turn shape-to-fit-text back on (restores <a:spAutoFit/> in the XML).
$sh16.TextFrame.AutoSize = 1

# Re-assert the original height in case re-enabling autofit nudged it by
# a fraction of a point; nudge the assigned value by tiny increments
# until the shape reports back the exact original height again.
$targetEmu16 = [math]::Round($origHeight16 * 12700)
$candidate16 = $origHeight16
for ($k = 0; $k -lt 1000; $k++) {
    $sh16.Height = $candidate16
    $resultEmu = [math]::Round($sh16.Height * 12700)
    if ($resultEmu -eq $targetEmu16) {
        break
    }
    $candidate16 = $candidate16 + 0.0000005
}
